$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix missing geographic data: underscore the city names in the shared strings
$ws.Range("A2").Value = "Huambo_City"
$ws.Range("A3").Value = "Luanda_City"

# Re-apply font to the data range, which causes Excel to mint a new cell style
$ws.Range("A2:I4").Style = "Normal"

# Move selection as recorded at save time
$ws.Range("I15:J15").Select()
